# New crime data collected — update the weekly CompStat report:
#  - bump the report "Volume/Number" and covering-week dates in the header
#  - refresh the weekly/28-day/YTD/2-year crime counts and % change figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 31   Number  21" -> "...22" ---
$ws.Range('A8').Value = 'Volume 31   Number  22'

# --- Header: reporting week "5/20/2024 .. 5/26/2024" -> "5/27/2024 .. 6/2/2024" ---
$ws.Range('C9').Value = 'Report Covering the Week  5/27/2024  Through  6/2/2024'

# --- Row 14 (Murder) ---
$ws.Range('L14').NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range('L14').Value = 100

# --- Row 15 (Rape) ---
$ws.Range('L15').Value = -66.666666666666
$ws.Range('M15').Value = -66.666666666666

# --- Row 16 (Robbery) ---
$ws.Range('D16').Value = 1
$ws.Range('E16').Value = 100
$ws.Range('F16').Value = 8
$ws.Range('G16').Value = 11
$ws.Range('H16').Value = -27.272727272727
$ws.Range('I16').Value = 52
$ws.Range('J16').Value = 44
$ws.Range('K16').Value = 18.181818181818
$ws.Range('L16').Value = 26.829268292682
$ws.Range('M16').Value = -38.823529411764
$ws.Range('N16').Value = -82.894736842105

# --- Row 17 (Fel. Assault) ---
$ws.Range('C17').Value = 5
$ws.Range('D17').Value = 5
$ws.Range('E17').Value = 0
$ws.Range('F17').Value = 19
$ws.Range('G17').Value = 23
$ws.Range('H17').Value = -17.391304347826
$ws.Range('I17').Value = 84
$ws.Range('J17').Value = 93
$ws.Range('K17').Value = -9.677419354838
$ws.Range('L17').Value = -15.151515151515
$ws.Range('M17').Value = 44.827586206896
$ws.Range('N17').Value = -70.526315789473

# --- Row 18 (Burglary) ---
$ws.Range('C18').Value = 3
$ws.Range('D18').NumberFormat = '#,##0'
$ws.Range('D18').Value = 1
$ws.Range('E18').NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range('E18').Value = 200
$ws.Range('F18').Value = 7
$ws.Range('G18').NumberFormat = '#,##0'
$ws.Range('G18').Value = 1
$ws.Range('H18').NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range('H18').Value = 600
$ws.Range('I18').Value = 36
$ws.Range('J18').Value = 27
$ws.Range('K18').Value = 33.333333333333
$ws.Range('L18').Value = -10
$ws.Range('M18').Value = 12.5
$ws.Range('N18').Value = -88.273615635179

# --- Row 19 (Gr. Larceny) ---
$ws.Range('C19').Value = 9
$ws.Range('D19').Value = 6
$ws.Range('E19').Value = 50
$ws.Range('F19').Value = 28
$ws.Range('G19').Value = 26
$ws.Range('H19').Value = 7.692307692307
$ws.Range('I19').Value = 148
$ws.Range('J19').Value = 119
$ws.Range('K19').Value = 24.369747899159
$ws.Range('L19').Value = -3.267973856209
$ws.Range('M19').Value = 146.666666666667
$ws.Range('N19').Value = 8.029197080291

# --- Row 20 (G.L.A.) ---
$ws.Range('C20').Value = 3
$ws.Range('F20').Value = 9
$ws.Range('G20').Value = 10
$ws.Range('H20').Value = -10
$ws.Range('I20').Value = 31
$ws.Range('K20').Value = -31.111111111111
$ws.Range('L20').Value = 19.230769230769
$ws.Range('M20').Value = 29.166666666666
$ws.Range('N20').Value = -75.590551181102

# --- Row 21 (TOTAL) ---
$ws.Range('C21').Value = 22
$ws.Range('D21').Value = 13
$ws.Range('E21').Value = 69.230769230769
$ws.Range('F21').Value = 73
$ws.Range('G21').Value = 71
$ws.Range('H21').Value = 2.816901408450
$ws.Range('I21').Value = 357
$ws.Range('J21').Value = 330
$ws.Range('K21').Value = 8.181818181818
$ws.Range('L21').Value = -4.032258064516
$ws.Range('M21').Value = 30.769230769230
$ws.Range('N21').Value = -70.520231213872

# --- Row 22 (Transit) ---
$ws.Range('G22').NumberFormat = '@'
$ws.Range('G22').Value = '0'
$ws.Range('H22').NumberFormat = '@'
$ws.Range('H22').Value = '***.*'

# --- Row 23 (Housing) ---
$ws.Range('G23').Value = 2
$ws.Range('J23').Value = 8
$ws.Range('K23').Value = -62.5

# --- Row 24 (Petit Larceny) ---
$ws.Range('C24').Value = 9
$ws.Range('D24').Value = 25
$ws.Range('E24').Value = -64
$ws.Range('F24').Value = 54
$ws.Range('G24').Value = 66
$ws.Range('H24').Value = -18.181818181818
$ws.Range('I24').Value = 290
$ws.Range('J24').Value = 298
$ws.Range('K24').Value = -2.684563758389
$ws.Range('L24').Value = -40.451745379876
$ws.Range('M24').Value = 124.806201550388

# --- Row 25 (Retail Theft) ---
$ws.Range('C25').Value = 2
$ws.Range('D25').Value = 8
$ws.Range('E25').Value = -75
$ws.Range('G25').Value = 21
$ws.Range('H25').Value = -28.571428571428
$ws.Range('I25').Value = 62
$ws.Range('J25').Value = 73
$ws.Range('K25').Value = -15.068493150684
$ws.Range('L25').Value = -71.296296296296

# --- Row 26 (Misd. Assault) ---
$ws.Range('C26').Value = 5
$ws.Range('D26').Value = 7
$ws.Range('E26').Value = -28.571428571428
$ws.Range('F26').Value = 32
$ws.Range('G26').Value = 25
$ws.Range('H26').Value = 28
$ws.Range('I26').Value = 140
$ws.Range('J26').Value = 122
$ws.Range('K26').Value = 14.754098360655
$ws.Range('L26').Value = -15.662650602409
$ws.Range('M26').Value = -28.571428571428

# --- Row 27 (UCR Rape*) ---
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = '0'
$ws.Range('L27').Value = -50

# --- Row 28 (Other Sex Crimes) ---
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '***.*'
$ws.Range('F28').Value = 1
$ws.Range('H28').Value = -50

# --- Row 29 (Shooting Vic.) ---
$ws.Range('G29').NumberFormat = '@'
$ws.Range('G29').Value = '0'
$ws.Range('H29').NumberFormat = '@'
$ws.Range('H29').Value = '***.*'
$ws.Range('N29').Value = -88.372093023255

# --- Row 30 (Shooting Inc.) ---
$ws.Range('G30').NumberFormat = '@'
$ws.Range('G30').Value = '0'
$ws.Range('H30').NumberFormat = '@'
$ws.Range('H30').Value = '***.*'
$ws.Range('N30').Value = -90.476190476190
